$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "b" (blue) color rows in column H (CLR) to "y" (yellow)
# for the extended color review. A new shared string "y" will be created.
$ws.Range("H4").Value = "y"
$ws.Range("H7").Value = "y"
$ws.Range("H10").Value = "y"
$ws.Range("H13").Value = "y"
$ws.Range("H16").Value = "y"
$ws.Range("H19").Value = "y"
$ws.Range("H22").Value = "y"

# Update the current selection to span the whole CLR column range
$ws.Range("H2:H22").Select()
